$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.523.52"
$ws.Range("E2").Value = "  -1.15%  "
# Row 3
$ws.Range("D3").Value = "1.912.93"
$ws.Range("E3").Value = "  -1.84%  "
# Row 4
$ws.Range("E4").Value = "  -0.06%  "
# Row 5
$ws.Range("D5").Value = "'239.60"
$ws.Range("E5").Value = "  -1.00%  "
# Row 6
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.16%  "
# Row 7
$ws.Range("D7").Value = "'0.4782"
$ws.Range("E7").Value = "  -1.67%  "
# Row 8
$ws.Range("D8").Value = "'0.2845"
$ws.Range("E8").Value = "  -2.76%  "
# Row 9
$ws.Range("D9").Value = "'0.06696"
$ws.Range("E9").Value = "  -2.16%  "
# Row 10
$ws.Range("D10").Value = "'18.68"
$ws.Range("E10").Value = "  -4.15%  "
# Row 11
$ws.Range("D11").Value = "'101.24"
$ws.Range("E11").Value = "  -3.85%  "
# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.916.14"
$ws.Range("E12").Value = "  -1.63%  "
# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07685"
$ws.Range("E13").Value = "  -0.82%  "
# Row 14
$ws.Range("D14").Value = "'5.226"
$ws.Range("E14").Value = "  -1.48%  "
# Row 15
$ws.Range("D15").Value = "'0.6699"
$ws.Range("E15").Value = "  -3.64%  "
# Row 16
$ws.Range("D16").Value = "30.525.06"
$ws.Range("E16").Value = "  -1.11%  "
# Row 17
$ws.Range("D17").Value = "'256.31"
$ws.Range("E17").Value = "  -7.04%  "
# Row 18
$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = "  -0.14%  "
# Row 19
$ws.Range("D19").Value = "'0.000007477"
$ws.Range("E19").Value = "  -3.04%  "
# Row 20
$ws.Range("D20").Value = "'12.67"
$ws.Range("E20").Value = "  -3.49%  "
# Row 21
$ws.Range("D21").Value = "'5.380"
$ws.Range("E21").Value = "  -1.14%  "
# Row 22
$ws.Range("E22").Value = "  -0.13%  "
# Row 23
$ws.Range("D23").Value = "'6.295"
$ws.Range("E23").Value = "  -2.62%  "
# Row 24
$ws.Range("D24").Value = "'9.337"
$ws.Range("E24").Value = "  -3.70%  "
# Row 25
$ws.Range("D25").Value = "'167.12"
$ws.Range("E25").Value = "  -0.19%  "
# Row 26
$ws.Range("D26").Value = "'19.09"
$ws.Range("E26").Value = "  -2.14%  "
# Row 27
$ws.Range("D27").Value = "'2.058"
$ws.Range("E27").Value = "  -4.68%  "
# Row 28
$ws.Range("D28").Value = "'4.764"
$ws.Range("E28").Value = "  +5.50%  "
# Row 29
$ws.Range("D29").Value = "'0.1009"
$ws.Range("E29").Value = "  -2.71%  "
# Row 30
$ws.Range("D30").Value = "'1.380"
$ws.Range("E30").Value = "  -0.96%  "
# Row 31
$ws.Range("E31").Value = "  -2.71%  "
# Row 32
$ws.Range("D32").Value = "'4.250"
$ws.Range("E32").Value = "  -2.62%  "
# Row 33
$ws.Range("D33").Value = "'0.04718"
$ws.Range("E33").Value = "  -2.59%  "
# Row 34
$ws.Range("D34").Value = "'0.7310"
$ws.Range("E34").Value = "  -1.91%  "
# Row 35
$ws.Range("D35").Value = "'1.109"
$ws.Range("E35").Value = "  -3.98%  "
# Row 36
$ws.Range("D36").Value = "'0.9997"
$ws.Range("E36").Value = "  -0.17%  "
# Row 37
$ws.Range("D37").Value = "'2.704"
$ws.Range("E37").Value = "  -0.97%  "
# Row 38
$ws.Range("D38").Value = "'0.01917"
$ws.Range("E38").Value = "  -3.38%  "
# Row 39
$ws.Range("D39").Value = "'2.612"
$ws.Range("E39").Value = "  -2.38%  "
# Row 40
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'6.231"
$ws.Range("E40").Value = "  -3.62%  "
# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'74.87"
$ws.Range("E41").Value = "  -2.34%  "
# Row 42
$ws.Range("D42").Value = "'1.971"
$ws.Range("E42").Value = "  -5.56%  "
# Row 43
$ws.Range("D43").Value = "'0.8620"
$ws.Range("E43").Value = "  -3.73%  "
# Row 44
$ws.Range("D44").Value = "'105.28"
$ws.Range("E44").Value = "  -2.47%  "
# Row 45
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.4243"
$ws.Range("E45").Value = "  -3.51%  "
# Row 46
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'0.9997"
$ws.Range("E46").Value = "  +0.02%  "
# Row 47
$ws.Range("D47").Value = "'7.386"
$ws.Range("E47").Value = "  -4.28%  "
# Row 48
$ws.Range("D48").Value = "'0.1200"
$ws.Range("E48").Value = "  -3.22%  "
# Row 49
$ws.Range("D49").Value = "'34.76"
$ws.Range("E49").Value = "  -2.28%  "
# Row 50
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "'907.29"
$ws.Range("E50").Value = "  -8.77%  "
# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'8.781"
$ws.Range("E51").Value = "  -4.10%  "
